# Auto-generated Excel COM-interop script to update cryptos list
# Commit: Updated cryptos list on Sat Jun 17 17:20:23 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '26.466.51'
$ws.Range("E2").Value = '  +2.03%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '1.728.59'
$ws.Range("E3").Value = '  +2.63%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  +0.17%  '

# Row 5: BNB
$ws.Range("D5").Value = '''244.40'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.62%  '

# Row 6: USDC
$ws.Range("E6").Value = '  +0.14%  '

# Row 7: XRP
$ws.Range("D7").Value = '''0.4810'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.48%  '

# Row 8: Cardano
$ws.Range("D8").Value = '''0.2681'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.75%  '

# Row 9: Dogecoin
$ws.Range("E9").Value = '  +1.26%  '

# Row 10: WrappedEther
$ws.Range("D10").Value = '1.732.74'
$ws.Range("E10").Value = '  +2.97%  '

# Row 11: TRON
$ws.Range("D11").Value = '''0.07132'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.59%  '

# Row 12: Solana
$ws.Range("D12").Value = '''15.77'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.75%  '

# Row 13: Polygon
$ws.Range("D13").Value = '''0.6186'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.42%  '

# Row 14: Polkadot
$ws.Range("D14").Value = '''4.556'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.80%  '

# Row 15: Litecoin
$ws.Range("E15").Value = '  +2.07%  '

# Row 16: Dai
$ws.Range("D16").Value = '''0.9999'
$ws.Range("D16").Style = "Normal"

# Row 17: WrappedBTC
$ws.Range("D17").Value = '26.479.73'
$ws.Range("E17").Value = '  +2.07%  '

# Row 18: BinanceUSD
$ws.Range("D18").Value = '''0.9998'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.08%  '

# Row 19: ShibaInu
$ws.Range("D19").Value = '''0.000006958'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.54%  '

# Row 20: Avalanche
$ws.Range("D20").Value = '''11.73'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.23%  '

# Row 21: WrappedliquidstakedEther2.0
$ws.Range("D21").Value = '1.954.13'
$ws.Range("E21").Value = '  +3.28%  '

# Row 22: Uniswap
$ws.Range("D22").Value = '''4.555'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.23%  '

# Row 23: Cosmos
$ws.Range("D23").Value = '''8.943'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.72%  '

# Row 24: Chainlink
$ws.Range("D24").Value = '''5.323'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.25%  '

# Row 25: Monero
$ws.Range("D25").Value = '''136.54'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.56%  '

# Row 26: EthereumClassic
$ws.Range("D26").Value = '''15.36'
$ws.Range("D26").Style = "Normal"

# Row 27: LidoDAOToken
$ws.Range("E27").Value = '  +4.61%  '

# Row 28: Toncoin
$ws.Range("D28").Value = '''1.405'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.98%  '

# Row 29: BitcoinCash
$ws.Range("E29").Value = '  +1.60%  '

# Row 30: InternetComputer(DFINITY)
$ws.Range("D30").Value = '''3.988'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.64%  '

# Row 31: Stellar
$ws.Range("D31").Value = '''0.08033'
$ws.Range("D31").Style = "Normal"

# Row 32: Filecoin
$ws.Range("D32").Value = '''3.743'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.94%  '

# Row 33: Hedera
$ws.Range("D33").Value = '''0.04565'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.52%  '

# Row 34: HuobiToken -> Frax
$ws.Range("B34").Value = 'Frax'
$ws.Range("C34").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D34").Value = '''0.9994'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.22%  '

# Row 35: ImmutableX -> HuobiToken
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '''2.615'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.75%  '

# Row 36: ARBITRUM -> ImmutableX
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '''0.6416'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.64%  '

# Row 37: TrustWalletToken -> ARBITRUM
$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D37").Value = '''0.9913'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.75%  '

# Row 38: RenderToken -> TrustWalletToken
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").Value = '''0.9439'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.06%  '

# Row 39: Quant -> RenderToken
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").Value = '''1.999'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +6.33%  '

# Row 40: MXToken -> Quant
$ws.Range("B40").Value = 'Quant'
$ws.Range("C40").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D40").Value = '''107.98'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.55%  '

# Row 41: PaxDollar -> MXToken
$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").Value = '''2.411'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.32%  '

# Row 42: VeChain -> PaxDollar
$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D42").Value = '''1.006'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.92%  '

# Row 43: FraxShare -> VeChain
$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").Value = '''0.01502'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.58%  '

# Row 44: TheSandbox -> FraxShare
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").Value = '''5.658'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +11.49%  '

# Row 45: Aptos -> TheSandbox
$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").Value = '''0.3922'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.97%  '

# Row 46: Algorand -> Aptos
$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").Value = '''7.020'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +13.68%  '

# Row 47: Cronos -> Algorand
$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").Value = '''0.1195'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +6.28%  '

# Row 48: Elrond -> Cronos
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '''0.05322'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.26%  '

# Row 49: EnergySwap -> Elrond
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Value = '''30.92'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.12%  '

# Row 50: NEARProtocol -> EnergySwap
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '''7.897'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.31%  '

# Row 51: Decentraland -> NEARProtocol
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = '''1.274'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.03%  '
